$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H31").Value = 474.75
$ws.Range("I31").Value = 299.66666
$ws.Range("K31").Value = 898.9999799999999
$ws.Range("M31").Value = -668.9999799999999
$ws.Range("H86").Value = 5266306.5
$ws.Range("I86").Value = 3865.1667
$ws.Range("J86").Value = 13159969
$ws.Range("K86").Value = 3865.1667
$ws.Range("L86").Value = 13159969
$ws.Range("M86").Value = -2742.1667
$ws.Range("N86").Value = -13162215
$ws.Range("H87").Value = 68542.10000000001
$ws.Range("J87").Value = 74953.336
$ws.Range("L87").Value = 74953.336
$ws.Range("N87").Value = -77449.336
$ws.Range("H88").Value = 7500.3335
$ws.Range("I88").Value = 8249.75
$ws.Range("J88").Value = 6001.5
$ws.Range("K88").Value = 8249.75
$ws.Range("L88").Value = 6001.5
$ws.Range("M88").Value = -7843.75
$ws.Range("N88").Value = -6813.5
$ws.Range("H89").Value = 5266306.5
$ws.Range("I89").Value = 3865.1667
$ws.Range("J89").Value = 13159969
$ws.Range("K89").Value = 19325.8335
$ws.Range("L89").Value = 65799845
$ws.Range("M89").Value = -13709.8335
$ws.Range("N89").Value = -65811077
$ws.Range("H90").Value = 68542.10000000001
$ws.Range("J90").Value = 74953.336
$ws.Range("L90").Value = 224860.008
$ws.Range("N90").Value = -237340.008
$ws.Range("H91").Value = 7500.3335
$ws.Range("I91").Value = 8249.75
$ws.Range("J91").Value = 6001.5
$ws.Range("K91").Value = 8249.75
$ws.Range("L91").Value = 6001.5
$ws.Range("M91").Value = -6845.75
$ws.Range("N91").Value = -8809.5
$ws.Range("H116").Value = 0
$ws.Range("I116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 0
$ws.Range("L116").ClearContents()
$ws.Range("M116").ClearContents()
$ws.Range("N116").Value = 0
$ws.Range("H129").Value = 15629.0625
$ws.Range("I129").Value = 23138
$ws.Range("J129").Value = 12215.909
$ws.Range("K129").Value = 69414
$ws.Range("L129").Value = 36647.727
$ws.Range("M129").Value = -64414
$ws.Range("N129").Value = -46647.727
$ws.Range("H137").Value = 3282.1555
$ws.Range("I137").Value = 1463.0834
$ws.Range("J137").Value = 5361.095
$ws.Range("K137").Value = 4389.2502
$ws.Range("L137").Value = 16083.285
$ws.Range("M137").Value = -1839.2502
$ws.Range("N137").Value = -21183.285
$ws.Range("H138").Value = 4391.875
$ws.Range("J138").Value = 6150.727
$ws.Range("L138").Value = 18452.181
$ws.Range("N138").Value = -28732.181

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3720.0833
$ws.Range("I32").Value = 3045.574
$ws.Range("J32").Value = 9790.666999999999
$ws.Range("K32").Value = 3045.574
$ws.Range("L32").Value = 9790.666999999999
$ws.Range("M32").Value = -2758.574
$ws.Range("N32").Value = -10364.667
$ws.Range("H63").Value = 3076
$ws.Range("I63").Value = 2913.5
$ws.Range("K63").Value = 2913.5
$ws.Range("M63").Value = -2227.5
$ws.Range("H66").Value = 3076
$ws.Range("I66").Value = 2913.5
$ws.Range("K66").Value = 14567.5
$ws.Range("M66").Value = -11135.5
$ws.Range("H88").Value = 1479.15
$ws.Range("I88").Value = 1597.5454
$ws.Range("J88").Value = 1334.4445
$ws.Range("K88").Value = 1597.5454
$ws.Range("L88").Value = 1334.4445
$ws.Range("M88").Value = -1191.5454
$ws.Range("N88").Value = -2146.4445
$ws.Range("H91").Value = 1479.15
$ws.Range("I91").Value = 1597.5454
$ws.Range("J91").Value = 1334.4445
$ws.Range("K91").Value = 1597.5454
$ws.Range("L91").Value = 1334.4445
$ws.Range("M91").Value = -193.5454
$ws.Range("N91").Value = -4142.4445

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 82966.336
$ws.Range("J82").Value = 99450
$ws.Range("L82").Value = 99450
$ws.Range("N82").Value = -100216
$ws.Range("H85").Value = 82966.336
$ws.Range("J85").Value = 99450
$ws.Range("L85").Value = 99450
$ws.Range("N85").Value = -102102
$ws.Range("H86").Value = 851591.75
$ws.Range("I86").Value = 896264.6
$ws.Range("K86").Value = 896264.6
$ws.Range("M86").Value = -895141.6
$ws.Range("H89").Value = 851591.75
$ws.Range("I89").Value = 896264.6
$ws.Range("K89").Value = 4481323
$ws.Range("M89").Value = -4475707
$ws.Range("H105").Value = 34901.91
$ws.Range("I105").Value = 59858.055
$ws.Range("K105").Value = 59858.055
$ws.Range("M105").Value = -58111.055
$ws.Range("H107").Value = 2310.2
$ws.Range("I107").Value = 2281.077
$ws.Range("K107").Value = 2281.077
$ws.Range("M107").Value = -361.0770000000002

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 289322.56
$ws.Range("I58").Value = 716020.0600000001
$ws.Range("J58").Value = 4857.5713
$ws.Range("K58").Value = 716020.0600000001
$ws.Range("L58").Value = 4857.5713
$ws.Range("M58").Value = -715817.0600000001
$ws.Range("N58").Value = -5263.5713
$ws.Range("H99").Value = 5209.8237
$ws.Range("J99").Value = 6779.3335
$ws.Range("L99").Value = 6779.3335
$ws.Range("N99").Value = -9775.333500000001
$ws.Range("H122").Value = 1993.0625
$ws.Range("I122").Value = 1420.6364
$ws.Range("J122").Value = 3252.4
$ws.Range("K122").Value = 4261.9092
$ws.Range("L122").Value = 9757.200000000001
$ws.Range("M122").Value = -1811.9092
$ws.Range("N122").Value = -14657.2
$ws.Range("H126").Value = 5209.8237
$ws.Range("J126").Value = 6779.3335
$ws.Range("L126").Value = 20338.0005
$ws.Range("N126").Value = -25278.0005
$ws.Range("H132").Value = 3906.0908
$ws.Range("I132").Value = 2788.238
$ws.Range("K132").Value = 8364.714
$ws.Range("M132").Value = -5834.714
$ws.Range("H136").Value = 289322.56
$ws.Range("I136").Value = 716020.0600000001
$ws.Range("J136").Value = 4857.5713
$ws.Range("K136").Value = 2148060.18
$ws.Range("L136").Value = 14572.7139
$ws.Range("M136").Value = -2145510.18
$ws.Range("N136").Value = -19672.7139

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 43709840
$ws.Range("I4").Value = 4352935.5
$ws.Range("K4").Value = 13058806.5
$ws.Range("M4").Value = -13058694.5
$ws.Range("H7").Value = 1150
$ws.Range("I7").Value = 300
$ws.Range("J7").Value = 2000
$ws.Range("K7").Value = 900
$ws.Range("L7").Value = 6000
$ws.Range("M7").Value = -788
$ws.Range("N7").Value = -6224
$ws.Range("H37").Value = 121107.16
$ws.Range("J37").Value = 121107.16
$ws.Range("L37").Value = 363321.48
$ws.Range("N37").Value = -363545.48
$ws.Range("H92").Value = 1668929.4
$ws.Range("I92").Value = 3335125.8
$ws.Range("K92").Value = 10005377.4
$ws.Range("M92").Value = -10004129.4

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 40004896
$ws.Range("I70").Value = 3940.9167
$ws.Range("J70").Value = 76928856
$ws.Range("K70").Value = 3940.9167
$ws.Range("L70").Value = 76928856
$ws.Range("M70").Value = -3670.9167
$ws.Range("N70").Value = -76929396
$ws.Range("H73").Value = 40004896
$ws.Range("I73").Value = 3940.9167
$ws.Range("J73").Value = 76928856
$ws.Range("K73").Value = 3940.9167
$ws.Range("L73").Value = 76928856
$ws.Range("M73").Value = -3004.9167
$ws.Range("N73").Value = -76930728
$ws.Range("H102").Value = 1245.65
$ws.Range("I102").Value = 1328.8
$ws.Range("J102").Value = 996.2
$ws.Range("K102").Value = 1328.8
$ws.Range("L102").Value = 996.2
$ws.Range("M102").Value = 293.2
$ws.Range("N102").Value = -4240.2
$ws.Range("H117").Value = 45000
$ws.Range("J117").Value = 45000
$ws.Range("L117").Value = 45000
$ws.Range("N117").Value = -51884
$ws.Range("H122").Value = 9136.75
$ws.Range("I122").Value = 8339.467000000001
$ws.Range("J122").Value = 10465.556
$ws.Range("K122").Value = 25018.401
$ws.Range("L122").Value = 31396.668
$ws.Range("M122").Value = -22568.401
$ws.Range("N122").Value = -36296.66800000001

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6599.905
$ws.Range("I7").Value = 7257.5
$ws.Range("J7").Value = 5723.1113
$ws.Range("K7").Value = 7257.5
$ws.Range("L7").Value = 5723.1113
$ws.Range("M7").Value = -7145.5
$ws.Range("N7").Value = -5947.1113
$ws.Range("H40").Value = 4194.6875
$ws.Range("I40").Value = 3729.7144
$ws.Range("J40").Value = 5082.364
$ws.Range("K40").Value = 3729.7144
$ws.Range("L40").Value = 5082.364
$ws.Range("M40").Value = -3593.7144
$ws.Range("N40").Value = -5354.364
$ws.Range("H115").Value = 70077.8
$ws.Range("J115").Value = 70077.8
$ws.Range("L115").Value = 70077.8
$ws.Range("N115").Value = -72427.8
$ws.Range("H126").Value = 6599.905
$ws.Range("I126").Value = 7257.5
$ws.Range("J126").Value = 5723.1113
$ws.Range("K126").Value = 21772.5
$ws.Range("L126").Value = 17169.3339
$ws.Range("M126").Value = -19302.5
$ws.Range("N126").Value = -22109.3339

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H56").Value = 50156.5
$ws.Range("I56").Value = 0
$ws.Range("J56").Value = 50156.5
$ws.Range("K56").Value = 0
$ws.Range("L56").ClearContents()
$ws.Range("M56").Value = 50156.5
$ws.Range("N56").Value = -51584.5
$ws.Range("H62").Value = 507500
$ws.Range("I62").Value = 1000000
$ws.Range("J62").Value = 15000
$ws.Range("K62").Value = 1000000
$ws.Range("L62").Value = 15000
$ws.Range("M62").Value = -999376
$ws.Range("N62").Value = -16248
$ws.Range("H65").Value = 507500
$ws.Range("I65").Value = 1000000
$ws.Range("J65").Value = 15000
$ws.Range("K65").Value = 5000000
$ws.Range("L65").Value = 75000
$ws.Range("M65").Value = -4996880
$ws.Range("N65").Value = -81240
$ws.Range("H81").Value = 7560.8887
$ws.Range("I81").Value = 1077.1333
$ws.Range("J81").Value = 39979.668
$ws.Range("K81").Value = 2154.2666
$ws.Range("L81").Value = 79959.336
$ws.Range("M81").Value = -1093.2666
$ws.Range("N81").Value = -82081.336
$ws.Range("H84").Value = 7560.8887
$ws.Range("I84").Value = 1077.1333
$ws.Range("J84").Value = 39979.668
$ws.Range("K84").Value = 10771.333
$ws.Range("L84").Value = 399796.68
$ws.Range("M84").Value = -5467.332999999999
$ws.Range("N84").Value = -410404.68
$ws.Range("H107").Value = 203180.4
$ws.Range("I107").Value = 203180.4
$ws.Range("K107").Value = 609541.2
$ws.Range("M107").Value = -607621.2
$ws.Range("H122").Value = 35719010
$ws.Range("I122").Value = 55559300
$ws.Range("K122").Value = 166677900
$ws.Range("M122").Value = -166675450
